$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The automated per-jornada loader now leaves "Capitan" (column AP) blank for
# every player row instead of writing an explicit "No"/"Si" flag.
$ws.Range("AP2:AP21").Value = ""

# Normalize "Suplente" (column F) typo: "Si" -> "Sí" (with accent) for the
# substitute players introduced later in the sheet (rows 13-21).
for ($r = 13; $r -le 21; $r++) {
    $ws.Range("F$r").Value = "Sí"
}
